$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("N15").Style = "Bad"
$ws.Range("N15").HorizontalAlignment = "Center"
